$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4
$ws.Range("A5").Value = 5
$ws.Range("A6").Value = 6
$ws.Range("A7").Value = 7

$ws.Range("A8").Select()
